# Applies the "poligonal" data-correction edit:
#  - adjusts four cumulative-distance readings in column C
#  - updates the active selection to reflect the user's last-used range

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C3").Value = 2364551
$ws.Range("C5").Value = 2813643
$ws.Range("C7").Value = 3063246
$ws.Range("C9").Value = 750440

$ws.Range("C2:C9").Select()
